$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B "Valor" shifts to C)
$ws.Columns.Item(2).Insert()

# Header row
$ws.Range("B1").Value = "Variável"
$ws.Range("C1").Value = "Valor"
$ws.Range("D1").Value = "Colocação"

# New "Variável" column for data rows 2-8
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 2).Value = "Diferença 2021-2012"
}

# "Colocação" ranking column, rows 2-7 only (Brasil row 8 left blank)
$rank = @("1º", "2º", "3º", "4º", "5º", "6º")
for ($i = 0; $i -lt $rank.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 4).Value = $rank[$i]
}
